# Insert a new daily price record for "Mango" / Vega Modelo de Temuco
# at row 163, pushing the existing rows 163-280 down to 164-281.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 163 (shifts 163..280 -> 164..281).
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new record's values.
$ws.Cells.Item(163, 1).Value  = 10
$ws.Cells.Item(163, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(163, 3).Value  = "La Araucanía"
$ws.Cells.Item(163, 4).Value  = 44574
$ws.Cells.Item(163, 5).Value  = 9
$ws.Cells.Item(163, 6).Value  = "Fruta"
$ws.Cells.Item(163, 7).Value  = 100108
$ws.Cells.Item(163, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(163, 9).Value  = 100108002
$ws.Cells.Item(163, 10).Value = "Mango"
$ws.Cells.Item(163, 11).Value = "Sin especificar"
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 1500
$ws.Cells.Item(163, 14).Value = 7000
$ws.Cells.Item(163, 15).Value = 7000
$ws.Cells.Item(163, 16).Value = 7000
$ws.Cells.Item(163, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(163, 18).Value = "Perú"
$ws.Cells.Item(163, 19).Value = 1750
$ws.Cells.Item(163, 20).Value = 4
